$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "567.80", "1.00", "0.0000245") that
# must stay text. A bare `.Value = "567.80"` lets Excel auto-coerce the
# literal to a number (567.8), dropping the formatting. Force the Text
# number format before the write, then restore the default "Normal" style
# afterwards so the saved cell has no stray style index (matches the
# original workbook, where these cells carry no explicit style).
function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "65.778.71"
$ws.Range("E2").Value = "  +1.25%  "
Set-TextValue "D3" "2.953.41"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.23%  "
Set-TextValue "D5" "567.80"
$ws.Range("E5").Value = "  -2.11%  "
Set-TextValue "D6" "159.91"
$ws.Range("E6").Value = "  +4.73%  "
$ws.Range("E7").Value = "  +0.10%  "
Set-TextValue "D8" "0.517"
$ws.Range("E8").Value = "  +0.82%  "
Set-TextValue "D9" "2.947.31"
$ws.Range("E9").Value = "  -1.06%  "
Set-TextValue "D10" "6.74"
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("E11").Value = "  +0.04%  "
Set-TextValue "D12" "0.457"
$ws.Range("E12").Value = "  +2.12%  "
Set-TextValue "D13" "0.0000245"
$ws.Range("E13").Value = "  +3.01%  "
Set-TextValue "D14" "34.20"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("E15").Value = "  -0.60%  "
Set-TextValue "D16" "65.839.09"
$ws.Range("E16").Value = "  +1.56%  "
Set-TextValue "D17" "3.444.96"
$ws.Range("E17").Value = "  -0.81%  "
Set-TextValue "D18" "6.92"
$ws.Range("E18").Value = "  +0.63%  "
Set-TextValue "D19" "2.953.63"
$ws.Range("E19").Value = "  -1.11%  "
Set-TextValue "D20" "446.27"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("E21").Value = "  +0.84%  "
Set-TextValue "D22" "0.677"
$ws.Range("E22").Value = "  -0.14%  "
Set-TextValue "D23" "7.19"
$ws.Range("E23").Value = "  -1.12%  "
Set-TextValue "D24" "82.48"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("E25").Value = "  +0.54%  "
Set-TextValue "D26" "12.18"
$ws.Range("E26").Value = "  -0.63%  "
Set-TextValue "D27" "10.02"
$ws.Range("E27").Value = "  -6.31%  "
$ws.Range("E28").Value = "  -0.02%  "
Set-TextValue "D29" "8.01"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  -0.05%  "
Set-TextValue "D32" "0.0₃0977"
$ws.Range("E32").Value = "  -4.63%  "
Set-TextValue "D33" "27.27"
$ws.Range("E33").Value = "  +2.48%  "
Set-TextValue "D34" "0.110"
$ws.Range("E34").Value = "  +0.19%  "
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  +0.18%  "
Set-TextValue "D36" "0.973"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("E38").Value = "  +0.57%  "
Set-TextValue "D39" "1.99"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("E40").Value = "  +1.89%  "
Set-TextValue "D41" "43.33"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -0.35%  "
Set-TextValue "D43" "2.80"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("E44").Value = "  +0.18%  "
Set-TextValue "D45" "385.58"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("E46").Value = "  +1.85%  "
Set-TextValue "D47" "2.724.14"
$ws.Range("E47").Value = "  -1.74%  "
Set-TextValue "D48" "130.67"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D50" "2.15"
$ws.Range("E50").Value = "  +5.34%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D51" "0.106"
$ws.Range("E51").Value = "  +1.19%  "
